# MLP and LSTM model, updated webapp2
#
# Adds two new data-source rows (Eurostat/Eurostag, EEX Transparency
# Platform) to the "production history" table, appends a clarifying
# remark ("existence of these data confirmed in [MAR2018]") to a
# description cell on the "weather history" sheet, and updates which
# sheet/cell is active/selected on several sheets.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# 1) "production history" sheet: add two rows to the table (Table3)
# ---------------------------------------------------------------
$wsProd = $wb.Worksheets.Item("production history")
$lo = $wsProd.ListObjects.Item(1)

$lo.ListRows.Add() | Out-Null
$lo.ListRows.Add() | Out-Null

# Fill in this particular order so newly-created shared strings line up
# the same way they did in the authored workbook.
$wsProd.Range("B5").Value = "https://www.eex-transparency.com/power/de/production/usage"
$wsProd.Range("D5").Value = "`"Power Production per Unit`", not accessible with bulk download or API, no additional power plant information, not even location, only for a very small past time period"
$wsProd.Range("A5").Value = "EEX Transparency Platform"

$wsProd.Range("A4").Value = "Eurostag"
$wsProd.Range("D4").Value = "not at WPP level, only at country lvel; monthly aggregated"
$wsProd.Range("B4").Value = "https://ec.europa.eu/eurostat/de/data/database"

$wsProd.Range("C4").Value = "free"
$wsProd.Range("C5").Value = "free"
$wsProd.Range("E4").Value = "Europe"
$wsProd.Range("E5").Value = "Europe"

# Give the new Link cells the same "Hyperlink" look as the existing ones
$wsProd.Range("B4").Style = "Hyperlink"
$wsProd.Range("B4").WrapText = $true
$wsProd.Range("B5").Style = "Hyperlink"
$wsProd.Range("B5").WrapText = $true

# ---------------------------------------------------------------
# 2) "weather history" sheet: extend the sensor-data remark
# ---------------------------------------------------------------
$wsWeatherHist = $wb.Worksheets.Item("weather history")
$wsWeatherHist.Range("E6").Value = "for yaw, pitch, stall control, WPPs measure wind speed and direction --> difficult ot obtain, not required, because good weather data available, not useful, because arbitrary locations should be selectable, formerly sensors behind nacelles --> lower wind speed, existence of these data confirmed in [MAR2018]"

# ---------------------------------------------------------------
# 3) Update active cell / selection on each sheet
# ---------------------------------------------------------------
$wsWPPs = $wb.Worksheets.Item("WPPs")
$wsWPPs.Activate()
$wsWPPs.Range("A5").Select()

$wsWeatherHist.Activate()
$wsWeatherHist.Range("B6").Select()

$wsProd.Activate()
$wsProd.Range("A2").Select()

$wsForecast = $wb.Worksheets.Item("weather forecast")
$wsForecast.Activate()
$wsForecast.Range("B5").Select()
